# Actualización automática 2025-11-11 10:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D34").Value = 2577.79
$ws1.Range("M34").Value = 551.98
$ws1.Range("L44").Value = 443.44
$ws1.Range("L60").Value = "4 de 58"
$ws1.Range("M60").Value = "4 de 58"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F34").Value = 2722.86
$ws2.Range("F44").Value = 443.44
$ws2.Range("F60").Value = 6694.04

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 3966.14
$ws3.Range("E3").Value = 10859.27
$ws3.Range("F3").Value = 0.2675231241496863

$ws3.Range("D11").Value = 1360.94
$ws3.Range("E11").Value = 14787.06
$ws3.Range("F11").Value = 0.08427916769878623

$ws3.Range("D12").Value = 803.5700000000001
$ws3.Range("E12").Value = 49503.43
$ws3.Range("F12").Value = 0.01597332379191763

$ws3.Range("D14").Value = 6694.04
$ws3.Range("E14").Value = 91167.84766749099
$ws3.Range("F14").Value = 0.06840293151450942

Write-Host "Edits applied."
